$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) column cells to retain exact text representation
$ws.Range("D2:D5").NumberFormat = "@"
$ws.Range("D9:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.524.40"
$ws.Range("D3").Value = "1.854.75"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").Value = "234.01"
$ws.Range("D9").Value = "0.06336"
$ws.Range("D10").Value = "17.66"
$ws.Range("D11").Value = "1.841.91"
$ws.Range("D12").Value = "0.07411"
$ws.Range("D13").Value = "5.026"
$ws.Range("D14").Value = "84.77"
$ws.Range("D15").Value = "0.6288"
$ws.Range("D16").Value = "30.468.22"
$ws.Range("D17").Value = "243.18"
$ws.Range("D18").Value = "1.000"
$ws.Range("D19").Value = "12.70"
$ws.Range("D20").Value = "0.000007351"
$ws.Range("D21").Value = "0.9992"
$ws.Range("D22").Value = "4.966"
$ws.Range("D23").Value = "5.980"
$ws.Range("D24").Value = "9.259"
$ws.Range("D25").Value = "162.51"
$ws.Range("D26").Value = "18.10"
$ws.Range("D27").Value = "1.890"
$ws.Range("D28").Value = "0.1014"
$ws.Range("D29").Value = "1.374"
$ws.Range("D30").Value = "4.048"
$ws.Range("D31").Value = "3.860"
$ws.Range("D32").Value = "0.04890"
$ws.Range("D33").Value = "1.141"
$ws.Range("D34").Value = "0.7058"
$ws.Range("D35").Value = "2.706"
$ws.Range("D36").Value = "0.01903"
$ws.Range("D37").Value = "2.685"
$ws.Range("D38").Value = "0.8730"
$ws.Range("D39").Value = "1.983"
$ws.Range("D40").Value = "105.47"
$ws.Range("D41").Value = "1.000"
$ws.Range("D42").Value = "0.4079"
$ws.Range("D43").Value = "5.499"
$ws.Range("D44").Value = "7.252"
$ws.Range("D45").Value = "62.71"
$ws.Range("D46").Value = "0.1203"
$ws.Range("D47").Value = "33.44"
$ws.Range("D48").Value = "8.532"
$ws.Range("D49").Value = "0.05530"
$ws.Range("D50").Value = "1.373"
$ws.Range("D51").Value = "0.3694"

$ws.Range("D2:D5").NumberFormat = "General"
$ws.Range("D2:D5").Style = "Normal"
$ws.Range("D9:D51").NumberFormat = "General"
$ws.Range("D9:D51").Style = "Normal"

# Update Volume(1h) (E) column text values
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("E8").Value = "  -0.21%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  +7.96%  "
$ws.Range("E11").Value = "  -1.13%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("E15").Value = "  -1.15%  "
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +4.89%  "
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("E33").Value = "  -1.75%  "
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  +1.37%  "
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("E43").Value = "  -1.69%  "
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -1.81%  "
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("E51").Value = "  -0.64%  "
